$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.269.56'
$ws.Range('E2').Value = '  -0.36%  '
# Row 3
$ws.Range('D3').Value = '3.566.75'
$ws.Range('E3').Value = '  +2.90%  '
# Row 4
$ws.Range('E4').Value = '  -0.15%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.79%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.23%  '
# Row 7
$ws.Range('D7').Value = '3.563.73'
$ws.Range('E7').Value = '  +2.87%  '
# Row 8
$ws.Range('E8').Value = '  +0.06%  '
# Row 9
$ws.Range('E9').Value = '  +2.30%  '
# Row 10
$ws.Range('E10').Value = '  -1.63%  '
# Row 11
$ws.Range('E11').Value = '  +1.52%  '
# Row 12
$ws.Range('E12').Value = '  -0.33%  '
# Row 13
$ws.Range('D13').Value = '4.175.26'
$ws.Range('E13').Value = '  +2.98%  '
# Row 14
$ws.Range('E14').Value = '  -0.55%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.30%  '
# Row 16
$ws.Range('D16').Value = '3.565.88'
$ws.Range('E16').Value = '  +2.91%  '
# Row 17
$ws.Range('D17').Value = '66.345.76'
$ws.Range('E17').Value = '  -0.44%  '
# Row 18
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.115'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.26%  '
# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.60%  '
# Row 20
$ws.Range('E20').Value = '  +0.13%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.28%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.68%  '
# Row 23
$ws.Range('E23').Value = '  +2.10%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.22%  '
# Row 25
$ws.Range('D25').Value = '3.710.89'
$ws.Range('E25').Value = '  +2.95%  '
# Row 26
$ws.Range('E26').Value = '  +0.04%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000122'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.51%  '
# Row 28
$ws.Range('E28').Value = '  +1.35%  '
# Row 29
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.45%  '
# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.52%  '
# Row 31
$ws.Range('E31').Value = '  -0.10%  '
# Row 32
$ws.Range('E32').Value = '  -2.05%  '
# Row 33
$ws.Range('E33').Value = '  -2.98%  '
# Row 34
$ws.Range('D34').Value = '3.566.12'
$ws.Range('E34').Value = '  +2.78%  '
# Row 35
$ws.Range('E35').Value = '  +1.36%  '
# Row 36
$ws.Range('E36').Value = '  -0.06%  '
# Row 37
$ws.Range('E37').Value = '  +1.12%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.90'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.77%  '
# Row 39
$ws.Range('E39').Value = '  +0.25%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.17%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '171.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.57%  '
# Row 42
$ws.Range('E42').Value = '  -2.65%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.18%  '
# Row 44
$ws.Range('E44').Value = '  +1.32%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.93'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.93%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.06%  '
# Row 47
$ws.Range('E47').Value = '  +2.96%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.04%  '
# Row 49
$ws.Range('E49').Value = '  +2.82%  '
# Row 50
$ws.Range('E50').Value = '  -1.36%  '
# Row 51
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.950'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.12%  '
